$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-5.25%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "48.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.20%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.184"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.12%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07754"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-4.97%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.504"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.36%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.338"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "13.99%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.551"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.93%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1223"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.82%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.53%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04677"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.82%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09371"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.35%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.01%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001260"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-4.90%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04178"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.88%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005845"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.76%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.328"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.06%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.243"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-7.95%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.74%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.028"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.86%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-5.65%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3040"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.37%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001276"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.09%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004092"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001353"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.37%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.94%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02569"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-7.55%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05800"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.37%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01077"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "71.14%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007938"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.25%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1419"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.97%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008340"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "8.54%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007691"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.64%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3065"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.92%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006991"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.29%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.36%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05669"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-7.55%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.26%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.36%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.36%"
